# Update countries & provincias Spain
# Applies the per-row data refresh described by the commit:
#  - Kazajistan (row 33): new D/E values
#  - Nepal / Kirguistan (rows 61-62): Kirguistan overtakes Nepal in the
#    ranking, so the two rows swap country + stats
#  - El Salvador (row 74): new D/E/G/H values
#  - Birmania / Principado de Andorra / Niger (rows 160-162): Birmania
#    overtakes both, shifting Andorra and Niger down one rank
#  - A1 "Datos actualizados" timestamp bumped to 07:48

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 33: Kazajistan ---------------------------------------------------
$ws.Cells.Item(33, 4).Value = 99018
$ws.Cells.Item(33, 5).Value = 5619

# --- Rows 61-62: Kirguistan / Nepal swap ranking --------------------------
$ws.Cells.Item(61, 1).Value = "Kirguistan"
$ws.Cells.Item(61, 2).Value = 44293
$ws.Cells.Item(61, 3).Value = 94
$ws.Cells.Item(61, 4).Value = 39599
$ws.Cells.Item(61, 5).Value = 3634
$ws.Cells.Item(61, 6).Value = 0
$ws.Cells.Item(61, 7).Value = 0
$ws.Cells.Item(61, 8).Value = 1060

$ws.Cells.Item(62, 1).Value = "Nepal"
$ws.Cells.Item(62, 2).Value = 44236
$ws.Cells.Item(62, 3).Value = 0
$ws.Cells.Item(62, 4).Value = 25561
$ws.Cells.Item(62, 5).Value = 18404
$ws.Cells.Item(62, 6).Value = 0
$ws.Cells.Item(62, 7).Value = 0
$ws.Cells.Item(62, 8).Value = 271

# --- Row 74: El Salvador ---------------------------------------------------
$ws.Cells.Item(74, 4).Value = 15369
$ws.Cells.Item(74, 5).Value = 9978
$ws.Cells.Item(74, 7).Value = 8
$ws.Cells.Item(74, 8).Value = 752

# --- Rows 160-162: Birmania / Principado de Andorra / Niger ---------------
$ws.Cells.Item(160, 1).Value = "Birmania"
$ws.Cells.Item(160, 2).Value = 1216
$ws.Cells.Item(160, 3).Value = 45
$ws.Cells.Item(160, 4).Value = 359
$ws.Cells.Item(160, 5).Value = 850
$ws.Cells.Item(160, 6).Value = 0
$ws.Cells.Item(160, 7).Value = 0
$ws.Cells.Item(160, 8).Value = 7

$ws.Cells.Item(161, 1).Value = "Principado de Andorra"
$ws.Cells.Item(161, 2).Value = 1215
$ws.Cells.Item(161, 3).Value = 0
$ws.Cells.Item(161, 4).Value = 928
$ws.Cells.Item(161, 5).Value = 234
$ws.Cells.Item(161, 6).Value = 0
$ws.Cells.Item(161, 7).Value = 0
$ws.Cells.Item(161, 8).Value = 53

$ws.Cells.Item(162, 1).Value = "Niger"
$ws.Cells.Item(162, 2).Value = 1177
$ws.Cells.Item(162, 3).Value = 0
$ws.Cells.Item(162, 4).Value = 1091
$ws.Cells.Item(162, 5).Value = 17
$ws.Cells.Item(162, 6).Value = 0
$ws.Cells.Item(162, 7).Value = 0
$ws.Cells.Item(162, 8).Value = 69

# --- A1: refresh timestamp --------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 5 de Septiembre de 2020 a las 07:48"
